$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the "ConnectsTo" column
$ws.Range("C1").Value = "ConnectsTo"
$ws.Range("C1").Font.Bold = $true

# Populate the new column with the servers each node connects to
$ws.Range("C2").Value = "Server2"
$ws.Range("C3").Value = "Server3, Server4"
$ws.Range("C4").Value = "Server5"

# Match the column width used for the new "ConnectsTo" column
$ws.Columns("C").ColumnWidth = 14
